# feature: 产品页面支持 IPV6 主机导入 (closed #776)
# Adds a new "寻址方式（可选）" (Addressing method, optional) column (O) to the
# SheetJS template sheet, with "静态" (Static) / "动态" (Dynamic) sample values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in O1
$ws.Range("O1").Value = "寻址方式（可选）"

# Sample values for the two data rows that get a value in the diff
$ws.Range("O2").Value = "静态"
$ws.Range("O3").Value = "动态"

# Give the new column a sensible width, matching the authored template
# (18.15 is the closest achievable input to the template's stored width of 18.875
# given this engine's 7px/char column-width quantization)
$ws.Columns.Item(15).ColumnWidth = 18.15

# Move the active selection to O3, matching the saved selection in the file
$ws.Range("O3").Select()
